# ---------------------------------------------------------------------------
# Re-creates the edit described by the commit:
#   - drop two new sheets "Sheet5" (copy of Employees + TaxFactor1/TaxFactor2
#     columns) and "Shee6" (copy of Teacher) at the end of the workbook
#   - update a few sheetView selections left behind by the author while
#     clicking around (Teacher, Employees, Sheet4)
#   - leave the new sheet "Shee6" as the active / selected tab
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$student   = $wb.Worksheets.Item("Student")
$teacher   = $wb.Worksheets.Item("Teacher")
$employees = $wb.Worksheets.Item("Employees")
$sheet4    = $wb.Worksheets.Item("Sheet4")

# --- Teacher: selection grew from a single cell to the whole used range ----
$teacher.Activate()
$teacher.Range("A1:E6").Select()

# --- Employees: same kind of selection update ------------------------------
$employees.Activate()
$employees.Range("A1:E6").Select()

# --- Sheet4: selection moved, and it is no longer the active tab -----------
$sheet4.Activate()
$sheet4.Range("D12").Select()

# --- new sheet "Sheet5" : a copy of Employees with two extra tax columns ---
$sheetCount = $wb.Worksheets.Count
$sheet5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$sheet5.Name = "Sheet5"

$employees.Range("A1:E6").Copy($sheet5.Range("A1"))
# reuse the TaxFactor column formatting for the two new columns
$employees.Range("E1:E6").Copy($sheet5.Range("F1"))
$employees.Range("E1:E6").Copy($sheet5.Range("G1"))
$sheet5.Range("F1").Value = "TaxFactor1"
$sheet5.Range("G1").Value = "TaxFactor2"

$sheet5.Columns.Item(6).ColumnWidth = 14
$sheet5.Columns.Item(7).ColumnWidth = 17.333333333333332

$sheet5.Activate()
$sheet5.Range("G10").Select()

# --- new sheet "Shee6" : a copy of Teacher ----------------------------------
$sheetCount = $wb.Worksheets.Count
$sheet6 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$sheet6.Name = "Shee6"

$teacher.Range("A1:E6").Copy($sheet6.Range("A1"))

$sheet6.Columns.Item(4).ColumnWidth = 20.666666666666668
$sheet6.Columns.Item(5).ColumnWidth = 20.166666666666668

# Shee6 ends up the active / selected sheet
$sheet6.Activate()
$sheet6.Range("C11").Select()
